# Fill in the test-case results for the population-projection worksheet.
# Rows 6-10 correspond to the five scenarios already labeled in column A
# (Regular US Data, Population Increase - High Birth Rate, Population
# Increase - High Migration, Population Decrease - High Death Rate,
# Population Low Birth Rate and Low Migration). Columns B:H hold the
# inputs/outputs for each scenario; column I is left blank, same as before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - Regular US Data
$ws.Range("B6").Value = 8
$ws.Range("C6").Value = 12
$ws.Range("D6").Value = 126
$ws.Range("E6").Value = 333100360
$ws.Range("F6").Value = 5
$ws.Range("G6").Value = 7821428
$ws.Range("H6").Value = 340921788

# Row 7 - Population Increase - High Birth Rate
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = 12
$ws.Range("D7").Value = 126
$ws.Range("E7").Value = 333100360
$ws.Range("F7").Value = 5
$ws.Range("G7").Value = 19647428
$ws.Range("H7").Value = 352747788

# Row 8 - Population Increase - High Migration
$ws.Range("B8").Value = 8
$ws.Range("C8").Value = 12
$ws.Range("D8").Value = 80
$ws.Range("E8").Value = 333100360
$ws.Range("F8").Value = 5
$ws.Range("G8").Value = 8541000
$ws.Range("H8").Value = 341641360

# Row 9 - Population Decrease - High Death Rate (G9 left blank, like source)
$ws.Range("B9").Value = 8
$ws.Range("C9").Value = 8
$ws.Range("D9").Value = 126
$ws.Range("E9").Value = 333100360
$ws.Range("F9").Value = 5
$ws.Range("H9").Value = 334351788

# Row 10 - Population Low Birth Rate and Low Migration (G10 left blank, like source)
$ws.Range("B10").Value = 14
$ws.Range("C10").Value = 12
$ws.Range("D10").Value = 200
$ws.Range("E10").Value = 333100360
$ws.Range("F10").Value = 5
$ws.Range("H10").Value = 332011617

# Leave the selection on G10, matching the saved view.
[void]$ws.Range("G10").Select()
